$d = $word.ActiveDocument

# 1. Replace the header text with the [Naglowek] placeholder.
$d.Content.Find.Execute("XV Ogólnopolski Konkurs Modeli Kartonowych o Puchar Wójta Gminy Przeciszów", `
    $true, $false, $false, $false, $false, $true, 1, $false, "[Naglowek]", 2)

# 2. Move the `_GoBack` bookmark from the (now-empty) third paragraph to the
#    end of the first paragraph (right after the "[Naglowek]" run, before
#    the paragraph mark).
$p1 = $d.Paragraphs(1)
$insertPos = $p1.Range.End - 1   # position right after the run's text, before the pilcrow

$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# Insert a temporary placeholder character so the target position is no
# longer the exact "end of paragraph text" boundary (a boundary that the
# bookmark-insertion logic mishandles), add the bookmark before it, then
# remove the placeholder again.
$tmpRng = $d.Range($insertPos, $insertPos)
$tmpRng.InsertAfter("X")

$bmRng = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

$delRng = $d.Range($insertPos, $insertPos + 1)
$delRng.Delete()

Write-Host "done"
